$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for each changed cell (old -> new per commit diff)
$newValues = [ordered]@{
    "E2" = "1"
    "F2" = "0.3333333333333333"
    "G2" = "0.002837"
    "H2" = "0.008510999999999999"
    "I2" = "7.108247730492929E-05"
    "J2" = "7.108247730492929E-05"
    "M2" = "1.097462"
    "N2" = "3.292386"
    "O2" = "0.3941457306284283"
    "P2" = "0.3941457306284283"
    "Q2" = "0.003113499693999999"
    "R2" = "0.028021497246"
    "S2" = "2.801685495223003E-05"
    "T2" = "2.801685495223002E-05"
    "E3" = "1"
    "F3" = "0.3333333333333333"
    "G3" = "0.002837"
    "H3" = "0.008510999999999999"
    "I3" = "7.108247730492929E-05"
    "J3" = "7.108247730492929E-05"
    "O3" = "0.244178053493144"
    "P3" = "0.2441780534931439"
    "Q3" = "0.001928850767"
    "R3" = "0.017359656903"
    "S3" = "1.735678094578822E-05"
    "T3" = "1.735678094578822E-05"
    "E4" = "1"
    "F4" = "0.3333333333333333"
    "G4" = "0.002837"
    "H4" = "0.008510999999999999"
    "I4" = "7.108247730492929E-05"
    "J4" = "7.108247730492929E-05"
    "M4" = "0.5872626666666667"
    "N4" = "1.761788"
    "O4" = "0.2109112414134909"
    "P4" = "0.2109112414134909"
    "Q4" = "0.001666064185333333"
    "R4" = "0.014994577668"
    "S4" = "1.499209353112893E-05"
    "T4" = "1.499209353112893E-05"
    "E5" = "1"
    "F5" = "0.3333333333333333"
    "G5" = "0.002837"
    "H5" = "0.008510999999999999"
    "I5" = "7.108247730492929E-05"
    "J5" = "7.108247730492929E-05"
    "M5" = "0.419791"
    "N5" = "1.259373"
    "O5" = "0.1507649744649369"
    "P5" = "0.1507649744649369"
    "Q5" = "0.001190947067"
    "R5" = "0.010718523603"
    "S5" = "1.071674787578212E-05"
    "T5" = "1.071674787578212E-05"
    "I6" = "0.3776915775490952"
    "J6" = "0.3776915775490952"
    "M6" = "1.097462"
    "N6" = "3.292386"
    "O6" = "0.3941457306284283"
    "P6" = "0.3941457306284283"
    "Q6" = "16.54335436398667"
    "R6" = "148.89018927588"
    "S6" = "0.1488655227852918"
    "T6" = "0.1488655227852918"
    "I7" = "0.3776915775490952"
    "J7" = "0.3776915775490952"
    "O7" = "0.244178053493144"
    "P7" = "0.2441780534931439"
    "S7" = "0.0922239942266929"
    "T7" = "0.09222399422669289"
    "I8" = "0.3776915775490952"
    "J8" = "0.3776915775490952"
    "M8" = "0.5872626666666667"
    "N8" = "1.761788"
    "O8" = "0.2109112414134909"
    "P8" = "0.2109112414134909"
    "Q8" = "8.852510974782223"
    "R8" = "79.67259877304001"
    "S8" = "0.07965939949229943"
    "T8" = "0.07965939949229943"
    "I9" = "0.3776915775490952"
    "J9" = "0.3776915775490952"
    "M9" = "0.419791"
    "N9" = "1.259373"
    "O9" = "0.1507649744649369"
    "P9" = "0.1507649744649369"
    "Q9" = "6.328010693593334"
    "R9" = "56.95209624234001"
    "S9" = "0.05694266104481107"
    "T9" = "0.05694266104481107"
    "G10" = "1.581618666666667"
    "H10" = "4.744856"
    "I10" = "0.03962825977384063"
    "J10" = "0.03962825977384063"
    "M10" = "1.097462"
    "N10" = "3.292386"
    "O10" = "0.3941457306284283"
    "P10" = "0.3941457306284283"
    "Q10" = "1.735766385157333"
    "R10" = "15.621897466416"
    "S10" = "0.01561930940209357"
    "T10" = "0.01561930940209357"
    "G11" = "1.581618666666667"
    "H11" = "4.744856"
    "I11" = "0.03962825977384063"
    "J11" = "0.03962825977384063"
    "O11" = "0.244178053493144"
    "P11" = "0.2441780534931439"
    "Q11" = "1.075328296898667"
    "R11" = "9.677954672088001"
    "S11" = "0.009676351334897063"
    "T11" = "0.009676351334897061"
    "G12" = "1.581618666666667"
    "H12" = "4.744856"
    "I12" = "0.03962825977384063"
    "J12" = "0.03962825977384063"
    "M12" = "0.5872626666666667"
    "N12" = "1.761788"
    "O12" = "0.2109112414134909"
    "P12" = "0.2109112414134909"
    "Q12" = "0.9288255958364445"
    "R12" = "8.359430362528002"
    "S12" = "0.00835804546395703"
    "T12" = "0.008358045463957029"
    "G13" = "1.581618666666667"
    "H13" = "4.744856"
    "I13" = "0.03962825977384063"
    "J13" = "0.03962825977384063"
    "M13" = "0.419791"
    "N13" = "1.259373"
    "O13" = "0.1507649744649369"
    "P13" = "0.1507649744649369"
    "Q13" = "0.6639492816986667"
    "R13" = "5.975543535288001"
    "S13" = "0.005974553572892968"
    "T13" = "0.005974553572892968"
    "G14" = "23.25273433333334"
    "H14" = "69.75820300000001"
    "I14" = "0.5826090801997593"
    "J14" = "0.5826090801997593"
    "M14" = "1.097462"
    "N14" = "3.292386"
    "O14" = "0.3941457306284283"
    "P14" = "0.3941457306284283"
    "Q14" = "25.51899232692867"
    "R14" = "229.670930942358"
    "S14" = "0.2296328815860907"
    "T14" = "0.2296328815860907"
    "G15" = "23.25273433333334"
    "H15" = "69.75820300000001"
    "I15" = "0.5826090801997593"
    "J15" = "0.5826090801997593"
    "O15" = "0.244178053493144"
    "P15" = "0.2441780534931439"
    "Q15" = "15.80932479862434"
    "R15" = "142.283923187619"
    "S15" = "0.1422603511506082"
    "T15" = "0.1422603511506082"
    "G16" = "23.25273433333334"
    "H16" = "69.75820300000001"
    "I16" = "0.5826090801997593"
    "J16" = "0.5826090801997593"
    "M16" = "0.5872626666666667"
    "N16" = "1.761788"
    "O16" = "0.2109112414134909"
    "P16" = "0.2109112414134909"
    "Q16" = "13.65546277188489"
    "R16" = "122.899164946964"
    "S16" = "0.1228788043637033"
    "T16" = "0.1228788043637033"
    "G17" = "23.25273433333334"
    "H17" = "69.75820300000001"
    "I17" = "0.5826090801997593"
    "J17" = "0.5826090801997593"
    "M17" = "0.419791"
    "N17" = "1.259373"
    "O17" = "0.1507649744649369"
    "P17" = "0.1507649744649369"
    "Q17" = "9.761288598524335"
    "R17" = "87.85159738671902"
    "S17" = "0.08783704309935708"
    "T17" = "0.08783704309935708"
}

foreach ($ref in $newValues.Keys) {
    $ws.Range($ref).Value = [double]$newValues[$ref]
}

Write-Host "Updated $($newValues.Count) cells"